# Auto-generated edit script: apply numeric corrections to Kujata_Profits sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5138.7
$ws.Range("I6").Value = 5138.7
$ws.Range("K6").Value = 15416.1
$ws.Range("M6").Value = -15304.1
$ws.Range("H64").Value = 3584.842
$ws.Range("I64").Value = 3483.2856
$ws.Range("J64").Value = 3869.2
$ws.Range("K64").Value = 3483.2856
$ws.Range("L64").Value = 3869.2
$ws.Range("M64").Value = -3235.2856
$ws.Range("N64").Value = -4365.2
$ws.Range("H67").Value = 3584.842
$ws.Range("I67").Value = 3483.2856
$ws.Range("J67").Value = 3869.2
$ws.Range("K67").Value = 3483.2856
$ws.Range("L67").Value = 3869.2
$ws.Range("M67").Value = -2625.2856
$ws.Range("N67").Value = -5585.2
$ws.Range("H70").Value = 1225
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 766.6667
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 2300.0001
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -2840.0001
$ws.Range("H73").Value = 1225
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 766.6667
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 2300.0001
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -4172.0001
$ws.Range("H135").Value = 615
$ws.Range("I135").Value = 268.75
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 2418.75
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = 116.25
$ws.Range("N135").Value = -23070

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1077.2222
$ws.Range("I74").Value = 836
$ws.Range("J74").Value = 3007
$ws.Range("K74").Value = 836
$ws.Range("L74").Value = 3007
$ws.Range("M74").Value = 38
$ws.Range("N74").Value = -4755
$ws.Range("H77").Value = 1077.2222
$ws.Range("I77").Value = 836
$ws.Range("J77").Value = 3007
$ws.Range("K77").Value = 4180
$ws.Range("L77").Value = 15035
$ws.Range("M77").Value = 188
$ws.Range("N77").Value = -23771
$ws.Range("H110").Value = 2294.1538
$ws.Range("I110").Value = 1860
$ws.Range("J110").Value = 2487.111
$ws.Range("K110").Value = 1860
$ws.Range("L110").Value = 2487.111
$ws.Range("M110").Value = 185
$ws.Range("N110").Value = -6577.111
$ws.Range("H122").Value = 1953.6364
$ws.Range("I122").Value = 2059
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 6177
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -3727
$ws.Range("N122").Value = -7600
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3854.0435
$ws.Range("I86").Value = 4260.9414
$ws.Range("J86").Value = 2701.1667
$ws.Range("K86").Value = 4260.9414
$ws.Range("L86").Value = 2701.1667
$ws.Range("M86").Value = -3137.9414
$ws.Range("N86").Value = -4947.1667
$ws.Range("H89").Value = 3854.0435
$ws.Range("I89").Value = 4260.9414
$ws.Range("J89").Value = 2701.1667
$ws.Range("K89").Value = 21304.707
$ws.Range("L89").Value = 13505.8335
$ws.Range("M89").Value = -15688.707
$ws.Range("N89").Value = -24737.8335
$ws.Range("H105").Value = 111112984
$ws.Range("I105").Value = 142859070
$ws.Range("J105").Value = 1695.5
$ws.Range("K105").Value = 142859070
$ws.Range("L105").Value = 1695.5
$ws.Range("M105").Value = -142857323
$ws.Range("N105").Value = -5189.5
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125001230
$ws.Range("I16").Value = 250001000
$ws.Range("J16").Value = 1482
$ws.Range("K16").Value = 250001000
$ws.Range("L16").Value = 1482
$ws.Range("M16").Value = -250000713
$ws.Range("N16").Value = -2056
$ws.Range("H31").Value = 2321.963
$ws.Range("I31").Value = 1169.8
$ws.Range("K31").Value = 1169.8
$ws.Range("M31").Value = -874.8
$ws.Range("H34").Value = 2321.963
$ws.Range("I34").Value = 1169.8
$ws.Range("K34").Value = 1169.8
$ws.Range("M34").Value = -967.8
$ws.Range("H62").Value = 7695234.5
$ws.Range("I62").Value = 3044
$ws.Range("K62").Value = 3044
$ws.Range("M62").Value = -2420
$ws.Range("H65").Value = 7695234.5
$ws.Range("I65").Value = 3044
$ws.Range("K65").Value = 15220
$ws.Range("M65").Value = -12100
$ws.Range("H99").Value = 1631.2106
$ws.Range("I99").Value = 1684.3334
$ws.Range("J99").Value = 1583.4
$ws.Range("K99").Value = 1684.3334
$ws.Range("L99").Value = 1583.4
$ws.Range("M99").Value = -186.3334
$ws.Range("N99").Value = -4579.4
$ws.Range("H113").Value = 125001230
$ws.Range("I113").Value = 250001000
$ws.Range("J113").Value = 1482
$ws.Range("K113").Value = 250001000
$ws.Range("L113").Value = 1482
$ws.Range("M113").Value = -249998830
$ws.Range("N113").Value = -5822
$ws.Range("H122").Value = 949.2
$ws.Range("I122").Value = 888
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 2664
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -214
$ws.Range("N122").Value = -9400
$ws.Range("H126").Value = 1631.2106
$ws.Range("I126").Value = 1684.3334
$ws.Range("J126").Value = 1583.4
$ws.Range("K126").Value = 5053.0002
$ws.Range("L126").Value = 4750.200000000001
$ws.Range("M126").Value = -2583.0002
$ws.Range("N126").Value = -9690.200000000001
$ws.Range("H132").Value = 2605.0715
$ws.Range("I132").Value = 1847.3
$ws.Range("K132").Value = 5541.9
$ws.Range("M132").Value = -3011.9

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 469.1111
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 361
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 1083
$ws.Range("M7").Value = -1388
$ws.Range("N7").Value = -1307

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1810
$ws.Range("I61").Value = 1736
$ws.Range("J61").Value = 1933.3334
$ws.Range("K61").Value = 1736
$ws.Range("L61").Value = 1933.3334
$ws.Range("M61").Value = -1534
$ws.Range("N61").Value = -2337.3334
$ws.Range("H112").Value = 47665.668
$ws.Range("J112").Value = 47665.668
$ws.Range("L112").Value = 47665.668
$ws.Range("N112").Value = -50619.668
$ws.Range("H113").Value = 1810
$ws.Range("I113").Value = 1736
$ws.Range("J113").Value = 1933.3334
$ws.Range("K113").Value = 1736
$ws.Range("L113").Value = 1933.3334
$ws.Range("M113").Value = 434
$ws.Range("N113").Value = -6273.3334
$ws.Range("H132").Value = 60397.707
$ws.Range("I132").Value = 1054.1428
$ws.Range("K132").Value = 3162.4284
$ws.Range("M132").Value = -632.4284000000002

